$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.510352849960327
$ws.Range("B1").Value = 2.77774715423584
$ws.Range("C1").Value = 2.423213958740234
$ws.Range("D1").Value = 2.613227128982544
$ws.Range("E1").Value = 2.872695684432983
